# Project_#RPA_Challenge code and deployment updated
# Bump the packaged RoboticEnterpriseFramework version shown on the
# "Deployment sheet" worksheet, and leave the workbook positioned/selected
# the way the author last left it before re-uploading to SharePoint.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Package details -> Package Name version bump: 1.0.1 -> 1.0.2
$ws.Range("C17").Value = "1.0.2"

# Leave the view scrolled/selected on the package-details row, matching
# where the author was working when the sheet was saved.
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C17").Select()
